$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "D" (Price) column sometimes holds plain-numeric-looking text
# (e.g. "0.998", "219.48"). Excel auto-detects such literals as numbers
# when assigned via .Value, which would silently normalise values like
# "0.590" -> 0.59 (dropping the significant trailing zero) and would also
# flip the cell from text to a numeric type. The source data keeps these
# as text, so force the Text number format before writing any D-column
# value that looks like a plain number.

$ws.Range('D2').Value = '30.340.38'
$ws.Range('E2').Value = '  +1.53%  '

$ws.Range('D3').Value = '1.677.75'
$ws.Range('E3').Value = '  +3.16%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.27%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.48'
$ws.Range('E5').Value = '  +2.48%  '

$ws.Range('E6').Value = '  +0.87%  '

$ws.Range('E7').Value = '  -0.24%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '29.83'
$ws.Range('E8').Value = '  +0.33%  '

$ws.Range('E9').Value = '  +2.13%  '

$ws.Range('E10').Value = '  +1.14%  '

$ws.Range('E11').Value = '  -1.19%  '

$ws.Range('D12').Value = '1.917.34'
$ws.Range('E12').Value = '  +3.09%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.73'
$ws.Range('E13').Value = '  +18.60%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.619'
$ws.Range('E14').Value = '  +8.66%  '

$ws.Range('D15').Value = '1.674.56'
$ws.Range('E15').Value = '  +2.94%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.01'
$ws.Range('E16').Value = '  +3.11%  '

$ws.Range('D17').Value = '30.333.11'
$ws.Range('E17').Value = '  +1.38%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.81'
$ws.Range('E18').Value = '  +1.38%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '248.56'
$ws.Range('E19').Value = '  +0.86%  '

$ws.Range('E20').Value = '  +2.17%  '

$ws.Range('E21').Value = '  -0.20%  '

$ws.Range('E22').Value = '  +4.59%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.11'
$ws.Range('E23').Value = '  +5.38%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.22'
$ws.Range('E24').Value = '  +4.93%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.04'
$ws.Range('E25').Value = '  +0.53%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.88'
$ws.Range('E26').Value = '  +1.27%  '

$ws.Range('E27').Value = '  +0.33%  '

$ws.Range('E28').Value = '  +2.73%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.30%  '

$ws.Range('E30').Value = '  +2.10%  '

$ws.Range('E31').Value = '  +0.86%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.49'
$ws.Range('E32').Value = '  +4.00%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.31'
$ws.Range('E33').Value = '  +3.67%  '

$ws.Range('D34').Value = '1.483.77'
$ws.Range('E34').Value = '  +3.84%  '

$ws.Range('E35').Value = '  +4.77%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.04'
$ws.Range('E36').Value = '  +0.02%  '

$ws.Range('E37').Value = '  +5.19%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '79.65'
$ws.Range('E38').Value = '  +13.49%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.590'
$ws.Range('E39').Value = '  +6.18%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.71'
$ws.Range('E40').Value = '  -5.77%  '

$ws.Range('E41').Value = '  +0.50%  '

$ws.Range('E42').Value = '  +3.37%  '

$ws.Range('E43').Value = '  +2.85%  '

$ws.Range('E44').Value = '  +1.66%  '

$ws.Range('E45').Value = '  -3.69%  '

$ws.Range('E46').Value = '  -0.24%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '52.88'
$ws.Range('E47').Value = '  -4.13%  '

$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '1.811.29'
$ws.Range('E48').Value = '  +2.47%  '

$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.44'
$ws.Range('E49').Value = '  -0.03%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '95.31'
$ws.Range('E50').Value = '  +6.39%  '

$ws.Range('E51').Value = '  +10.43%  '
